# Regenerate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 2923c61f-...-md file's zh-cn / de-de round trip after a new
# handback report was produced.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-17 06:42:32"

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-17 06:42:27"
$zhcn.Range("K2").Value = "2016-08-17 06:42:45"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-17 06:42:32"
$dede.Range("K2").Value = "2016-08-17 06:42:52"
